# Update the division problems in the table to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "267÷9="; new = "380÷6="},
    @{old = "976÷8="; new = "770÷2="},
    @{old = "208÷4="; new = "549÷3="},
    @{old = "941÷2="; new = "366÷8="},
    @{old = "189÷6="; new = "630÷9="},
    @{old = "551÷3="; new = "861÷6="},
    @{old = "954÷4="; new = "716÷5="},
    @{old = "191÷5="; new = "871÷8="},
    @{old = "789÷6="; new = "179÷9="},
    @{old = "793÷4="; new = "967÷6="},
    @{old = "616÷7="; new = "305÷7="},
    @{old = "567÷9="; new = "386÷4="},
    @{old = "530÷8="; new = "993÷7="},
    @{old = "489÷2="; new = "362÷5="},
    @{old = "705÷3="; new = "571÷2="},
    @{old = "439÷8="; new = "958÷2="},
    @{old = "245÷5="; new = "816÷8="},
    @{old = "714÷5="; new = "389÷6="},
    @{old = "458÷7="; new = "963÷2="},
    @{old = "650÷6="; new = "886÷9="},
    @{old = "932÷9="; new = "549÷3="},
    @{old = "681÷7="; new = "656÷7="},
    @{old = "238÷4="; new = "165÷5="},
    @{old = "788÷2="; new = "874÷2="},
    @{old = "968÷3="; new = "472÷3="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
